$wb = $excel.ActiveWorkbook

# "fa9b0075-b570-4301-a336-d2c339cecb02" finished a later localization pass;
# refresh the "Latest Handoff Datetime" for the e3558659 row that feeds the
# newly generated handback, in both locale sheets.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "2016-03-13 06:41:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value = "2016-03-13 06:41:49"
